$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("H12").Value = 0
